$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for row 134 (Papaya, Vega Modelo de
# Temuco). The previous contents of row 134 are kept as historical data by
# inserting a fresh row right after it (new row 135), pushing the old row 135
# down to row 136.

$ws.Rows("135:135").Insert()

# Row 135 (new) gets the values that used to live in row 134 before the edit.
$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 45173
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = "Fruta"
$ws.Range("G135").Value = 100108
$ws.Range("H135").Value = "Tropicales y subtropicales"
$ws.Range("I135").Value = 100108004
$ws.Range("J135").Value = "Papaya"
$ws.Range("K135").Value = "Cultivar IV Región"
$ws.Range("L135").Value = "Primera"
$ws.Range("M135").Value = 150
$ws.Range("N135").Value = 24000
$ws.Range("O135").Value = 24000
$ws.Range("P135").Value = 24000
$ws.Range("Q135").Value = "$/bandeja 10 kilos"
$ws.Range("R135").Value = "Provincia del Elquí"
$ws.Range("S135").Value = 2400
$ws.Range("T135").Value = 10

# Row 134 is updated in place with the new, current observation.
$ws.Range("D134").Value = 45239
$ws.Range("M134").Value = 65
$ws.Range("N134").Value = 2600
$ws.Range("O134").Value = 2600
$ws.Range("P134").Value = 2600
$ws.Range("Q134").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S134").Value = 2600
$ws.Range("T134").Value = 1
